# Skosmos instances workbook - "Better query parsing" edit
#
# Inserts a new "context" column (C) before the existing "timeout" column
# (which shifts right, C -> D, and the old D "note" column shifts to E).
# Row 6 (Legilux) gets the text "wildcard" in the new context column; every
# other data row gets 0 in it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- insert the new column ----------------------------------------------
# Shift C:D -> D:E and create a fresh, blank column C.
$ws.Columns("C:C").Insert()

# --- header ---------------------------------------------------------------
$ws.Range("C1").Value2 = "context"

# --- data rows --------------------------------------------------------------
# Default every row to 0 ...
$ws.Range("C2:C10").Value2 = 0
# ... except row 6 (Legilux), which records the "wildcard" note instead.
$ws.Range("C6").Value2 = "wildcard"

# --- column widths ----------------------------------------------------------
# Keep the existing look: narrow "context"/"timeout" columns, a wider note
# column for the (now longer) remark text.
$ws.Columns("C:C").ColumnWidth = 7.666666666666667   # ~8.57
$ws.Columns("D:D").ColumnWidth = 7.333333333333333   # ~8.14
$ws.Columns("E:E").ColumnWidth = 21                  # ~21.86

# --- selection / view state -------------------------------------------------
$ws.Range("F22").Select()
